$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need NumberFormat forced to Text
# so Excel stores them as strings (matching the original inlineStr cells),
# then the style is reset back to Normal to avoid introducing new styling.
$textCells = $ws.Range("D5,D6,D7,D10,D16,D19,D20,D22,D23,D24,D25,D29,D30,D34,D36,D37,D39,D40,D41,D42,D43,D44,D45,D46,D47,D49,D50")
# NumberFormat/Style on a multi-area Range only affects the first area in
# this COM implementation, so apply per-area explicitly.
foreach ($area in $textCells.Areas) {
    $area.NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '56.968.01'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '2.394.38'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '504.33'
$ws.Range('E5').Value = '  -1.77%  '
$ws.Range('D6').Value = '131.61'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '2.406.64'
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('D10').Value = '0.0963'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('E13').Value = '  -4.77%  '
$ws.Range('D14').Value = '2.824.71'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('D15').Value = '56.892.09'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '21.70'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').Value = '2.448.53'
$ws.Range('E18').Value = '  +4.17%  '
$ws.Range('D19').Value = '10.18'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '309.43'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').Value = '6.33'
$ws.Range('E22').Value = '  +3.81%  '
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '65.22'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -3.47%  '
$ws.Range('D29').Value = '7.47'
$ws.Range('E29').Value = '  +3.99%  '
$ws.Range('D30').Value = '171.19'
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('D31').Value = '0.0₃0722'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').Value = '5.88'
$ws.Range('E34').Value = '  -3.91%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '17.91'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('D39').Value = '3.82'
$ws.Range('E39').Value = '  +2.70%  '
$ws.Range('D40').Value = '36.59'
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('D41').Value = '0.800'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('D42').Value = '1.43'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = '130.62'
$ws.Range('E43').Value = '  +7.32%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '3.35'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '4.82'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.565'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '251.37'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').Value = '0.0486'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').Value = '16.96'
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('E51').Value = '  +0.68%  '

# Restore default styling on the cells we forced to text format
foreach ($area in $textCells.Areas) {
    $area.Style = "Normal"
}

Write-Host "Applied cryptos update."
